$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "tahun" column header is renamed to "tgl_mitra_diterima"
# (K1 holds the last header in the mitra import/seeder sheet).
$ws.Range("K1").Value = "tgl_mitra_diterima"
